$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: update the rpc-reply message-id UUID ---
$f2 = $ws.Range("F2").Value2
$f2 = $f2 -replace 'urn:uuid:2e9090e4-ee21-4d82-94be-084d8bba156b', 'urn:uuid:93d96bc0-ab76-4426-a06b-0dc75e61653e'
$ws.Range("F2").Value = $f2

# --- G2: update protocol identifier (add namespace prefix) and rename BGP_65000 -> default ---
$g2 = $ws.Range("G2").Value2
$g2 = $g2 -replace [regex]::Escape('<identifier>BGP</identifier>'), '<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>'
$g2 = $g2 -replace [regex]::Escape('<name>BGP_65000</name>'), '<name>default</name>'
$ws.Range("G2").Value = $g2
